# Insert a new data row at row 234 (pushing the former rows 234-286 down to
# 235-287, which is exactly what the target diff shows: every old row from
# 234 onward reappears one row lower with identical values, the used range
# grows from A1:R286 to A1:R287, and a brand-new record appears at row 234).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(234).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(234, 1).Value  = 9
$ws.Cells.Item(234, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(234, 3).Value  = "Metropolitana"
$ws.Cells.Item(234, 4).Value  = 44711
$ws.Cells.Item(234, 5).Value  = 13
$ws.Cells.Item(234, 6).Value  = 100112043
$ws.Cells.Item(234, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(234, 8).Value  = "Sin especificar"
$ws.Cells.Item(234, 9).Value  = "Primera"
$ws.Cells.Item(234, 10).Value = 52
$ws.Cells.Item(234, 11).Value = 18000
$ws.Cells.Item(234, 12).Value = 20000
$ws.Cells.Item(234, 13).Value = 19000
$ws.Cells.Item(234, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(234, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(234, 16).Value = 317
$ws.Cells.Item(234, 17).Value = 60
$ws.Cells.Item(234, 18).Value = "Hortaliza"
